$d = $word.ActiveDocument

function Set-ParagraphText($para, [string]$newText) {
    $r = $d.Range($para.Range.Start, $para.Range.End - 1)
    $r.Text = $newText
}

function Find-ParagraphContaining([string]$needle) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

# --- 1. Intro bio paragraph + following "Accepting..." paragraph ---
# Merge: drop "from management through to presentation." and the whole
# "Accepting contracted development work within New Zealand and Denmark."
# paragraph, leaving just "...data lifecycle." in the bio paragraph.
$introPara = Find-ParagraphContaining "Comfortable with all aspects of the data lifecycle"

# First remove the paragraph mark that ends the intro paragraph so the
# following "Accepting ..." paragraph folds into it (true paragraph
# deletion, not just a text-content clear).
$markRng = $d.Range($introPara.Range.End - 1, $introPara.Range.End)
$markRng.Delete()

# Now trim the merged paragraph's tail down to a single "." after
# "data lifecycle".
$introPara = Find-ParagraphContaining "Comfortable with all aspects of the data lifecycle"
$fullIntro = $introPara.Range.Text
$idx = $fullIntro.IndexOf("Comfortable with all aspects of the data lifecycle")
$startPos = $introPara.Range.Start + $idx
$endPos = $introPara.Range.End - 1
$mergeRange = $d.Range($startPos, $endPos)
$mergeRange.Text = "Comfortable with all aspects of the data lifecycle."

# --- 2. Skills table updates ---
# "JavaScript, CSS/SASS, Bootstrap 4, React" -> "TypeScript, CSS/SASS, Bootstrap, React"
$p = Find-ParagraphContaining "CSS/SASS"
Set-ParagraphText $p "TypeScript, CSS/SASS, Bootstrap, React"

# "Custom data-processing APIs & workflows" -> "Esri ArcPy & Geoprocessing Services"
$p = Find-ParagraphContaining "Custom data-processing APIs"
Set-ParagraphText $p "Esri ArcPy & Geoprocessing Services"

# "Esri ArcPy geoprocessing services" -> "Node.js & Express backend JavaScript APIs"
$p = Find-ParagraphContaining "Esri ArcPy geoprocessing services"
Set-ParagraphText $p "Node.js & Express backend JavaScript APIs"

# --- 3. Experience section: job title update ---
# "Contracted Senior Advisor to the Abley Ltd. software team"
#   -> "Contracted Senior Developer with the Abley Ltd. software team"
$p = Find-ParagraphContaining "Contracted Senior Advisor"
Set-ParagraphText $p "Contracted Senior Developer with the Abley Ltd. software team"

Write-Output "done"
